# Apply updated cryptocurrency price / 1h-volume values to the "cryptos" worksheet.
# Price values in column D are free-form text (e.g. using "." as a thousands
# separator), so cells whose new value could otherwise be auto-parsed by Excel
# as a number are forced to Text format first so they stay text, matching the
# original inline-string cell contents.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.621.91"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "3.446.17"
$ws.Range("E3").Value = "  +2.26%  "
$ws.Range("E4").Value = "  -0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "581.18"
$ws.Range("E5").Value = "  +1.53%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "149.22"
$ws.Range("E6").Value = "  +9.14%  "
$ws.Range("D7").Value = "3.446.69"
$ws.Range("E7").Value = "  +2.34%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("E10").Value = "  +1.73%  "
$ws.Range("E11").Value = "  +3.40%  "
$ws.Range("E12").Value = "  +1.61%  "
$ws.Range("D13").Value = "4.036.12"
$ws.Range("E13").Value = "  +2.28%  "
$ws.Range("E14").Value = "  +7.76%  "
$ws.Range("E15").Value = "  -0.46%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.0000176"
$ws.Range("E16").Value = "  +2.45%  "
$ws.Range("D17").Value = "3.445.58"
$ws.Range("E17").Value = "  +2.19%  "
$ws.Range("D18").Value = "61.728.54"
$ws.Range("E18").Value = "  +1.22%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.28"
$ws.Range("E19").Value = "  +8.36%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "14.36"
$ws.Range("E20").Value = "  +3.19%  "
$ws.Range("E21").Value = "  +1.14%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "388.91"
$ws.Range("E22").Value = "  +4.02%  "
$ws.Range("E23").Value = "  +2.76%  "
$ws.Range("D24").Value = "3.589.04"
$ws.Range("E24").Value = "  +2.20%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "73.07"
$ws.Range("E25").Value = "  +2.89%  "
$ws.Range("E26").Value = "  +0.25%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.999"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("E29").Value = "  +2.79%  "
$ws.Range("E30").Value = "  +3.80%  "
$ws.Range("E31").Value = "  +0.08%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.53"
$ws.Range("E32").Value = "  -13.64%  "
$ws.Range("E33").Value = "  +1.71%  "
$ws.Range("E34").Value = "  +1.52%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "24.05"
$ws.Range("E36").Value = "  +1.73%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "7.03"
$ws.Range("E37").Value = "  +2.52%  "
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("E39").Value = "  +1.53%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "166.10"
$ws.Range("E40").Value = "  +1.02%  "
$ws.Range("E41").Value = "  +3.26%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "27.11"
$ws.Range("E42").Value = "  +12.89%  "
$ws.Range("E43").Value = "  +2.20%  "
$ws.Range("E44").Value = "  +2.52%  "
$ws.Range("E45").Value = "  -0.01%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "42.29"
$ws.Range("E46").Value = "  +1.78%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.70"
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("D48").Value = "2.603.96"
$ws.Range("E48").Value = "  +6.10%  "
$ws.Range("E50").Value = "  +2.48%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "23.25"
$ws.Range("E51").Value = "  +1.57%  "
